# Weekly update: insert a new price record for "Albahaca" at
# Terminal La Palmera de La Serena, shifting the existing rows
# (122..188) down by one to (123..189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 122; everything below (old rows
# 122-188) slides down to 123-189, matching the diff exactly.
$ws.Rows("122:122").Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(122, 1).Value  = 8
$ws.Cells.Item(122, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(122, 3).Value  = "Coquimbo"
$ws.Cells.Item(122, 4).Value  = 45086
$ws.Cells.Item(122, 5).Value  = 4
$ws.Cells.Item(122, 6).Value  = 100112052
$ws.Cells.Item(122, 7).Value  = "Albahaca"
$ws.Cells.Item(122, 8).Value  = "Sin especificar"
$ws.Cells.Item(122, 9).Value  = "Primera"
$ws.Cells.Item(122, 10).Value = 960
$ws.Cells.Item(122, 11).Value = 3000
$ws.Cells.Item(122, 12).Value = 3500
$ws.Cells.Item(122, 13).Value = 3250
$ws.Cells.Item(122, 14).Value = "$/paquete"
$ws.Cells.Item(122, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(122, 16).Value = 3250
$ws.Cells.Item(122, 17).Value = 1
$ws.Cells.Item(122, 18).Value = "Hortaliza"
